$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.246.57"
$ws.Range("E2").Value = "  +6.53%  "
$ws.Range("D3").Value = "3.112.92"
$ws.Range("E3").Value = "  +4.48%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.28%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.106.43"
$ws.Range("E8").Value = "  +4.54%  "
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.152"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.69%  "
$ws.Range("E13").Value = "  +8.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.84%  "
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "3.628.24"
$ws.Range("E16").Value = "  +4.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").Value = "63.168.64"
$ws.Range("E18").Value = "  +6.37%  "
$ws.Range("D19").Value = "3.111.16"
$ws.Range("E19").Value = "  +4.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.97%  "
$ws.Range("E21").Value = "  +4.56%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  +7.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.67%  "
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.95%  "
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("E34").Value = "  +13.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +16.38%  "
$ws.Range("E36").Value = "  +6.45%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +21.47%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "438.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("D42").Value = "2.914.52"
$ws.Range("E42").Value = "  +6.20%  "
$ws.Range("E43").Value = "  +5.64%  "
$ws.Range("E44").Value = "  +11.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.32%  "
$ws.Range("E46").Value = "  +8.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.08%  "
